$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Mental" count was revised from N=43 to N=44; update the label cell.
$ws.Range("A5").Value = "Mental (N=44)"

# Widen column A so the longer labels are fully visible.
# (29.17 "characters" resolves to the stored OOXML width of exactly 30.)
$ws.Columns.Item(1).ColumnWidth = 29.17

# Leave the active selection on A5, matching where the edit was made.
$ws.Range("A5").Select()
